# Raul's Log - append new activity rows to the "Logs" sheet.
# Mirrors a block of rows that was typed/pasted in starting at row 1223
# (rows 1217-1222 are intentionally left blank/unused, same as the rest
# of the sheet), ending with the new dated entries for 1/9/2017 (serial 42744).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# --- three mostly-blank spacer rows (only B/C/D carry the table's usual
# borders; row 1224 also keeps a bordered F cell) -----------------------
$ws.Range("B3:D3").Copy() | Out-Null
$ws.Range("B1223:D1223").PasteSpecial(-4122) | Out-Null
$ws.Range("B1225:D1225").PasteSpecial(-4122) | Out-Null

$ws.Range("B1224:D1224").PasteSpecial(-4122) | Out-Null
$ws.Range("F7").Copy() | Out-Null
$ws.Range("F1224").PasteSpecial(-4122) | Out-Null

# --- the actual log entries --------------------------------------------
$entries = @(
    @{Row=1226; A="Demo";     C="1600"; D="ACE"; E="007"},
    @{Row=1227; A="Demo";     C="1600"; D="ACE"; E="009"},
    @{Row=1228; A="Demo";     C="1630"; D="OSG"; E="1008"},
    @{Row=1229; A="Demo";     C="1630"; D="OSG"; E="2008"},
    @{Row=1230; A="Demo";     C="1730"; D="OSG"; E="2009"},
    @{Row=1231; A="Demo";     C="1730"; D="HNE"; E="141"},
    @{Row=1232; A="Demo";     C="1730"; D="OSG"; E="2028"},
    @{Row=1233; A="Demo";     C="1800"; D="OSG"; E="1001"},
    @{Row=1234; A="Demo";     C="1830"; D="OSG"; E="2004"},
    @{Row=1235; A="Demo";     C="1900"; D="DB";  E="0004"},
    @{Row=1236; A="Demo";     C="1900"; D="DB";  E="0009"},
    @{Row=1237; A="Demo";     C="1900"; D="DB";  E="2116"},
    @{Row=1238; A="Demo";     C="1900"; D="HNE"; E="032"},
    @{Row=1239; A="Demo";     C="1900"; D="OSG"; E="2003"},
    @{Row=1240; A="Demo";     C="1630"; D="SSB"; E="W141"; F="PC, neck mic and podium mic,  audio cable for laptop!"},
    @{Row=1241; A="Operator"; C="1700"; D="SSB"; E="W141"; F="Operate event from 17:00-18:00"}
)

foreach ($e in $entries) {
    $r = $e.Row
    $ws.Cells.Item($r, 1).Value = $e.A
    $ws.Cells.Item($r, 2).Value = 42744
    $ws.Cells.Item($r, 3).Value = $e.C
    $ws.Cells.Item($r, 4).Value = $e.D
    $ws.Cells.Item($r, 5).Value = $e.E
    if ($e.F) {
        $ws.Cells.Item($r, 6).Value = $e.F
    }
}

$ws.Range("A1241").Select()
